# Fruta / hortaliza, semanal
# Insert a new weekly record at row 616 (Femacal de La Calera - Naranja, Valencia)
# pushing the existing rows 616:656 down to 617:657.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 616; this shifts rows
# 616-656 down to 617-657, carrying their formatting (incl. the date
# number format on column D) along with them.
$ws.Rows.Item(616).Insert()

# Populate the newly inserted row 616 with the new weekly price record.
$ws.Range("A616").Value = 3
$ws.Range("B616").Value = "Femacal de La Calera"
$ws.Range("C616").Value = "Coquimbo"
$ws.Range("D616").Value = 44610
$ws.Range("E616").Value = 5
$ws.Range("F616").Value = "Fruta"
$ws.Range("G616").Value = 100102
$ws.Range("H616").Value = "Cítricos"
$ws.Range("I616").Value = 100102005
$ws.Range("J616").Value = "Naranja"
$ws.Range("K616").Value = "Valencia"
$ws.Range("L616").Value = "Primera"
$ws.Range("M616").Value = 125
$ws.Range("N616").Value = 7000
$ws.Range("O616").Value = 7500
$ws.Range("P616").Value = 7260
$ws.Range("Q616").Value = "$/malla 13 kilos"
$ws.Range("R616").Value = "Provincia de Quillota"
$ws.Range("S616").Value = 558
$ws.Range("T616").Value = 13
